$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6262181997299194
$ws.Range("B1").Value = 1.590180397033691
$ws.Range("C1").Value = 5.004844188690186
$ws.Range("D1").Value = 1.633498668670654
$ws.Range("E1").Value = 0.9060161113739014
